$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.903.18'
$ws.Range("E2").Value = '  -1.95%  '

# Row 3
$ws.Range("D3").Value = '3.639.69'
$ws.Range("E3").Value = '  -0.10%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.56'
$ws.Range("E5").Value = '  -1.64%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.99'
$ws.Range("E6").Value = '  -3.00%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.634'
$ws.Range("E7").Value = '  +4.27%  '

# Row 8
$ws.Range("D8").Value = '3.632.76'
$ws.Range("E8").Value = '  +0.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.73'
$ws.Range("E11").Value = '  +14.88%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.623'
$ws.Range("E12").Value = '  +3.06%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '48.51'
$ws.Range("E13").Value = '  -3.11%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000284'
$ws.Range("E14").Value = '  -1.54%  '

# Row 15
$ws.Range("D15").Value = '4.225.49'
$ws.Range("E15").Value = '  -0.10%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '671.00'
$ws.Range("E16").Value = '  -4.11%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '9.07'
$ws.Range("E17").Value = '  +1.24%  '

# Row 18
$ws.Range("D18").Value = '3.634.70'
$ws.Range("E18").Value = '  -1.81%  '

# Row 19
$ws.Range("D19").Value = '70.913.64'

# Row 20
$ws.Range("E20").Value = '  -0.42%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.85'
$ws.Range("E21").Value = '  -3.43%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.50'
$ws.Range("E22").Value = '  -1.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.945'
$ws.Range("E23").Value = '  +1.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.24'
$ws.Range("E24").Value = '  -4.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.99'
$ws.Range("E25").Value = '  -3.74%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.92'
$ws.Range("E26").Value = '  -2.84%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.79'
$ws.Range("E27").Value = '  -2.36%  '

# Row 28
$ws.Range("E28").Value = '  +0.05%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.85'
$ws.Range("E29").Value = '  -1.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.68'
$ws.Range("E30").Value = '  -1.52%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.22'
$ws.Range("E31").Value = '  +0.88%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  -4.72%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.65'

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.39'
$ws.Range("E34").Value = '  -5.38%  '

# Row 35
$ws.Range("E35").Value = '  -4.61%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '575.53'
$ws.Range("E36").Value = '  -1.62%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.10'
$ws.Range("E37").Value = '  -1.71%  '

# Row 38
$ws.Range("E38").Value = '  -0.45%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '58.49'
$ws.Range("E39").Value = '  -2.22%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0454'
$ws.Range("E41").Value = '  -1.16%  '

# Row 42
$ws.Range("D42").Value = '3.550.61'
$ws.Range("E42").Value = '  -2.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.141'
$ws.Range("E43").Value = '  -2.34%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.44'
$ws.Range("E45").Value = '  -3.81%  '

# Row 46
$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D46").Value = '0.0₃0732'
$ws.Range("E46").Value = '  -5.75%  '

# Row 47
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.05'
$ws.Range("E47").Value = '  +6.49%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.68'
$ws.Range("E48").Value = '  -3.77%  '

# Row 49
$ws.Range("E49").Value = '  +1.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.20'
$ws.Range("E50").Value = '  +2.57%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.88'
$ws.Range("E51").Value = '  -4.24%  '
